$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.731.56'
$ws.Range('E2').Value = '  -4.02%  '
$ws.Range('D3').Value = '2.275.73'
$ws.Range('E3').Value = '  -5.50%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''542.61'
$ws.Range('E5').Value = '  -2.80%  '
$ws.Range('D6').Value = '''130.65'
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''0.567'
$ws.Range('E8').Value = '  -3.15%  '
$ws.Range('D9').Value = '2.275.52'
$ws.Range('E9').Value = '  -5.35%  '
$ws.Range('D10').Value = '''0.0997'
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('D11').Value = '''5.42'
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('E13').Value = '  -5.55%  '
$ws.Range('D14').Value = '''23.37'
$ws.Range('E14').Value = '  -5.37%  '
$ws.Range('D15').Value = '2.691.42'
$ws.Range('E15').Value = '  -5.15%  '
$ws.Range('D16').Value = '57.798.74'
$ws.Range('E16').Value = '  -3.76%  '
$ws.Range('E17').Value = '  -4.96%  '
$ws.Range('D18').Value = '2.297.53'
$ws.Range('E18').Value = '  -6.46%  '
$ws.Range('D19').Value = '''10.50'
$ws.Range('E19').Value = '  -6.33%  '
$ws.Range('D20').Value = '''4.23'
$ws.Range('E20').Value = '  -6.18%  '
$ws.Range('D21').Value = '''310.96'
$ws.Range('E21').Value = '  -4.66%  '
$ws.Range('D22').Value = '''6.38'
$ws.Range('E22').Value = '  -6.10%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '''62.64'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').Value = '''0.166'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -7.14%  '
$ws.Range('D28').Value = '''1.28'
$ws.Range('E28').Value = '  -7.95%  '
$ws.Range('D29').Value = '''1.72'
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('D30').Value = '''169.87'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E31').Value = '  -7.06%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '''5.69'
$ws.Range('E33').Value = '  -6.80%  '
$ws.Range('D34').Value = '''0.377'
$ws.Range('E34').Value = '  -6.28%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '''17.59'
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '''1.22'
$ws.Range('E38').Value = '  -8.32%  '
$ws.Range('E39').Value = '  -7.11%  '
$ws.Range('D40').Value = '''37.80'
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('E41').Value = '  -7.44%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''285.45'
$ws.Range('E42').Value = '  -12.06%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''138.75'
$ws.Range('E43').Value = '  -6.60%  '
$ws.Range('D44').Value = '''3.38'
$ws.Range('E44').Value = '  -5.10%  '
$ws.Range('D45').Value = '''0.0944'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').Value = '''0.0496'
$ws.Range('E46').Value = '  -3.75%  '
$ws.Range('D47').Value = '''0.548'
$ws.Range('E47').Value = '  -4.78%  '
$ws.Range('D48').Value = '''18.11'
$ws.Range('E48').Value = '  -8.98%  '
$ws.Range('D49').Value = '''0.0210'
$ws.Range('E49').Value = '  -5.01%  '
$ws.Range('D50').Value = '''10.96'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').Value = '''16.38'
$ws.Range('E51').Value = '  -4.64%  '
